# SCHEDULE_Fabián Becerra.xlsx - "Actualización valores Schedule, Task"
#
# - Turn off concurrent calculation for the workbook (calcPr concurrentCalc="0")
# - On Hoja1 (the active sheet):
#     * G14 goes from 3.2 to 3.5
#     * I14 goes from 3.2 to 3.5 (I15 = G15+I14 recalculates from 3.2 to 3.5)
#     * the current selection moves from H22 to I15

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# calcPr concurrentCalc="0"
$excel.MultiThreadedCalculation.Enabled = $false

# Updated task values
$ws.Range("G14").Value = 3.5
$ws.Range("I14").Value = 3.5

# New active cell/selection
$ws.Range("I15").Select()
